# Append a new row 24 to each of the 4 lifter log sheets, mirroring the
# existing row structure (time / length / id / actual-length / checksum
# plus their decimal counterparts).

$wb = $excel.ActiveWorkbook

$rows = @{
    "ROW35-FE-LIFTER" = @{
        A = 45736.28303427083
        B = "0x01,0x90"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c,"
        D = "0x01,0x82"
        E = "0xd"
        F = 400
        G = [double]"5.68631262647114e+23"
        H = 386
        I = 13
    }
    "ROW35-MID-LIFTER" = @{
        A = 45736.1332133912
        B = "0x01,0x90"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
        D = "0x01,0x82"
        E = "0xe"
        F = 400
        G = [double]"5.68631262647114e+23"
        H = 386
        I = 14
    }
    "ROW02-FE-LIFTER" = @{
        A = 45736.28025276621
        B = "0x01,0x90"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
        D = "0x01,0x82"
        E = "0x3"
        F = 400
        G = [double]"5.68631262647114e+23"
        H = 386
        I = 3
    }
    "ROW02-MID-LIFTER" = @{
        A = 45736.34179211806
        B = "0x01,0x90"
        C = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
        D = "0x01,0x82"
        E = "0x3"
        F = 400
        G = [double]"9.85046333984776e+23"
        H = 386
        I = 3
    }
}

foreach ($ws in $wb.Worksheets) {
    $data = $rows[$ws.Name]
    if ($data -eq $null) { continue }

    $newRow = 24

    # Column A: timestamp, numeric serial, same date/time style as the
    # rows above it (style index 2 / "YYYY-MM-DD HH:MM:SS").
    $cellA = $ws.Cells.Item($newRow, 1)
    $cellA.Value = $data.A
    $cellA.NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($newRow, 2).Value = $data.B
    $ws.Cells.Item($newRow, 3).Value = $data.C
    $ws.Cells.Item($newRow, 4).Value = $data.D
    $ws.Cells.Item($newRow, 5).Value = $data.E
    $ws.Cells.Item($newRow, 6).Value = $data.F
    $ws.Cells.Item($newRow, 7).Value = $data.G
    $ws.Cells.Item($newRow, 8).Value = $data.H
    $ws.Cells.Item($newRow, 9).Value = $data.I
}
